$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# number-format before assignment, otherwise Excel auto-converts "68.50" -> 68.5
# etc. and drops the trailing zero / fixed formatting. We restore the default
# "Normal" style afterwards so no stray formatting is left behind.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D14", "D19", "D23", "D24", "D25", "D28", "D30", "D32", "D36", "D38", "D40", "D42", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "49.002.10"
$ws.Range("E2").Value = "  -1.60%  "
$ws.Range("D3").Value = "2.624.27"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "111.23"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "322.41"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("D10").Value = "39.61"
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").Value = "19.69"
$ws.Range("E11").Value = "  -4.97%  "
$ws.Range("E12").Value = "  -1.80%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("D14").Value = "7.23"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "3.037.43"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "2.630.62"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "49.001.47"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "3.00"
$ws.Range("E19").Value = "  -3.47%  "
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "268.74"
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").Value = "68.50"
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "10.15"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "34.97"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -4.75%  "
$ws.Range("D32").Value = "49.37"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "18.98"
$ws.Range("E36").Value = "  -4.26%  "
$ws.Range("E37").Value = "  +4.28%  "
$ws.Range("D38").Value = "2.03"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "127.34"
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "22.17"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("E43").Value = "  -4.42%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "2.059.55"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +6.68%  "
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").Value = "  -5.34%  "
$ws.Range("E48").Value = "  -4.46%  "
$ws.Range("D49").Value = "8.89"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "5.18"
$ws.Range("E50").Value = "  -3.71%  "
$ws.Range("D51").Value = "58.57"
$ws.Range("E51").Value = "  +1.18%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
